$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fill in what used to be the trailing empty row (row 63): it now holds a
# real diary entry dated 19/05/2022 (serial 44700) - a 2 hour "Réalisation"
# entry about SQL-injection protection / exception handling.
$ws.Range("A63").Value = 44700
$ws.Range("B63").Value = "Réalisation"
$ws.Range("C63").Value = 2
$ws.Range("D63").Value = "Protection contre injections SQL dans l'ajout d'article et gestion d'Exceptions"

# The long description wraps onto two lines, same as the other multi-line rows.
$ws.Rows.Item(63).RowHeight = 30

# Add a new trailing, still-empty row 64 that only carries the date - this is
# the new "next entry" placeholder row (mirrors how row 63 looked before).
$ws.Range("A64").Value = 44700
$ws.Range("A63").Copy()
$ws.Range("A64").PasteSpecial(-4122)
$ws.Range("A64").Value = 44700

# Grow the table ("Tableau1") so its range / autofilter covers the new row.
$tbl = $ws.ListObjects.Item(1)
$tbl.Resize($ws.Range("A1:F64"))

# Leave the selection on the newly added blank cell, like the author did.
$ws.Range("B64").Select() | Out-Null

$excel.CutCopyMode = $false
